# murine-source.xlsx update:
#  - add "sex", "is_embryo", "date_of_birth_or_fertilization", "is_deceased",
#    "date_of_death", and "euthanization_method" columns (plus "strain" /
#    "strain_rrid" headers) ahead of the existing columns on the
#    "Export as TSV" sheet
#  - add two new lookup sheets: "sex list" and "euthanization_method list"
#  - add data validation + comments for the new columns

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Export as TSV")

# ---------------------------------------------------------------------------
# 1. Make room for 8 new leading columns (A:H). This shifts the existing
#    columns A:I (local_lifespan_data .. cage_enhancements) to I:Q, along with
#    their data validation rules.
# ---------------------------------------------------------------------------
$ws1.Range("A1:H1").EntireColumn.Insert()

# Copy the header style (bold, centered, wrap text) from the shifted
# local_lifespan_data header onto the new header cells.
$ws1.Range("I1").Copy()
$ws1.Range("A1:H1").PasteSpecial(-4122, -4142, $false, $false)

# ---------------------------------------------------------------------------
# 2. Fill in the new header values.
# ---------------------------------------------------------------------------
$ws1.Range("A1").Value = "strain"
$ws1.Range("B1").Value = "strain_rrid"
$ws1.Range("C1").Value = "sex"
$ws1.Range("D1").Value = "is_embryo"
$ws1.Range("E1").Value = "date_of_birth_or_fertilization"
$ws1.Range("F1").Value = "is_deceased"
$ws1.Range("G1").Value = "date_of_death"
$ws1.Range("H1").Value = "euthanization_method"

# ---------------------------------------------------------------------------
# 3. Add the two new lookup sheets, positioned right after "Export as TSV"
#    and before "room_health_status list":
#      sex list                  (index 2)
#      euthanization_method list (index 3)
# ---------------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("room_health_status list")
$wb.Worksheets.Add($beforeSheet) | Out-Null
$wb.Worksheets.Add($beforeSheet) | Out-Null

$wb.Worksheets.Item(2).Name = "sex list"
$wb.Worksheets.Item(3).Name = "euthanization_method list"

$sexSheet = $wb.Worksheets.Item("sex list")
$sexSheet.Range("A1").Value = "M"
$sexSheet.Range("A2").Value = "F"

$euthSheet = $wb.Worksheets.Item("euthanization_method list")
$euthSheet.Range("A1").Value = "Carbon dioxide asphixiation"
$euthSheet.Range("A2").Value = "Inhaled anesthetic agents"
$euthSheet.Range("A3").Value = "Injected anesthetic agents"
$euthSheet.Range("A4").Value = "Cervical dislocation"
$euthSheet.Range("A5").Value = "Decapitation"
$euthSheet.Range("A6").Value = "Hypothermia"
$euthSheet.Range("A7").Value = "Rapid freezing"
$euthSheet.Range("A8").Value = "Other"

# ---------------------------------------------------------------------------
# 4. Data validation for the new columns.
# ---------------------------------------------------------------------------

# sex (C) -> 'sex list'
$rngSex = $ws1.Range("C2:C1048576")
$rngSex.Validation.Add(3, 1, 1, "'sex list'!`$A`$1:`$A`$2")
$rngSex.Validation.ErrorTitle = "Value must come from list"
$rngSex.Validation.ErrorMessage = "Value must be one of: M / F."
$rngSex.Validation.IgnoreBlank = $true
$rngSex.Validation.InCellDropdown = $true
$rngSex.Validation.ShowInput = $true
$rngSex.Validation.ShowError = $true

# is_embryo (D) -> boolean
$rngEmbryo = $ws1.Range("D2:D1048576")
$rngEmbryo.Validation.Add(3, 1, 1, "TRUE,FALSE")
$rngEmbryo.Validation.ErrorTitle = "Not a boolean"
$rngEmbryo.Validation.ErrorMessage = "The values in this column must be `"TRUE`" or `"FALSE`"."
$rngEmbryo.Validation.IgnoreBlank = $true
$rngEmbryo.Validation.InCellDropdown = $true
$rngEmbryo.Validation.ShowInput = $true
$rngEmbryo.Validation.ShowError = $true

# is_deceased (F) -> boolean
$rngDeceased = $ws1.Range("F2:F1048576")
$rngDeceased.Validation.Add(3, 1, 1, "TRUE,FALSE")
$rngDeceased.Validation.ErrorTitle = "Not a boolean"
$rngDeceased.Validation.ErrorMessage = "The values in this column must be `"TRUE`" or `"FALSE`"."
$rngDeceased.Validation.IgnoreBlank = $true
$rngDeceased.Validation.InCellDropdown = $true
$rngDeceased.Validation.ShowInput = $true
$rngDeceased.Validation.ShowError = $true

# euthanization_method (H) -> 'euthanization_method list'
$rngEuth = $ws1.Range("H2:H1048576")
$rngEuth.Validation.Add(3, 1, 1, "'euthanization_method list'!`$A`$1:`$A`$8")
$rngEuth.Validation.ErrorTitle = "Value must come from list"
$rngEuth.Validation.ErrorMessage = "Value must come from euthanization_method list."
$rngEuth.Validation.IgnoreBlank = $true
$rngEuth.Validation.InCellDropdown = $true
$rngEuth.Validation.ShowInput = $true
$rngEuth.Validation.ShowError = $true

# ---------------------------------------------------------------------------
# 5. Comments. The column insert does not relocate existing cell comments, so
#    remove the stale ones (still anchored to the old A1:I1 addresses) and
#    re-create all 17 comments at their correct locations.
# ---------------------------------------------------------------------------
for ($c = 1; $c -le 9; $c++) {
    $cell = $ws1.Cells.Item(1, $c)
    if ($cell.Comment -ne $null) {
        $cell.Comment.Delete()
    }
}

$ws1.Range("A1").AddComment("Jackson Labs nomenclature. When mutant alleles are part of the strain name, use `"<`" and `">`" to indicate the superscripted alleles. For example, C57BL/6J-KitW-39J should be entered as `"C57BL/6J-Kit<W-39J>`", where `"W-39J`" would be the portion of the string displayed as superscripted text. For further information, see the `"Quick Guide to Mouse Nomenclature`" (https://resources.jax.org/guides/quick-guide-to-mouse-nomenclature).")
$ws1.Range("B1").AddComment("The Research Resource Identifier (RRID) (https://scicrunch.org/resources/data/source/nlx_154697-1/search) for the strain. An example is 'RRID:MGI:3713213'")
$ws1.Range("C1").AddComment("The sex of the mouse.")
$ws1.Range("D1").AddComment("Is the source an embryo? Use either 'True' or 'False'.")
$ws1.Range("E1").AddComment("The date when the mouse/embryo was born/fertilized. If the hours/minutes are not known, use '00:00'.")
$ws1.Range("F1").AddComment("Is the source deceased? Use either 'True' or 'False'.")
$ws1.Range("G1").AddComment("The date when the mouse/embryo died. If the hours/minutes are not known, use '00:00'.")
$ws1.Range("H1").AddComment("If the source was euthanized, select the method of euthanization.")
$ws1.Range("I1").AddComment("A free text description of how long mice live within the local environment. It is recommended to provide the median or maximum values for murine lifespans.")
$ws1.Range("J1").AddComment("A description of the pathogen and opportunist exclusion level of the room where the source is housed.")
$ws1.Range("K1").AddComment("The temperature value in Celsius of the room where the source is housed. An example is `"23`".")
$ws1.Range("L1").AddComment("The rack setup type in which the source is housed.")
$ws1.Range("M1").AddComment("The light cycle in the room where the source is housed. `"Standard/default`" refers to 12-hour photoperiods (e.g., lights on at 7:00 AM, lights off at 7:00 PM). `"Longer photoperiods`" refers to 14-hour photoperiods (e.g., lights on at 7:00 AM, lights off at 9:00 PM). `"Reverse lightcycles`" means that the the timing of the 12-hour photoperiod is reversed (.e.g, lights on at 7:00 PM, lights off at 7:00 AM).")
$ws1.Range("N1").AddComment("The type of cage bedding in the cage where the source is housed.")
$ws1.Range("O1").AddComment("A free text description of the source's diet.")
$ws1.Range("P1").AddComment("A free text description of the source's water supply, including any treatments to the water.")
$ws1.Range("Q1").AddComment("Environmental enrichments present in the source’s cage.")

Write-Host "Done updating murine-source workbook."
